$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.759.40"
$ws.Range("E2").Value = "  +4.66%  "
$ws.Range("D3").Value = "1.609.44"
$ws.Range("E3").Value = "  +3.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.45"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  +1.40%  "
$ws.Range("E6").Value = "  +6.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "26.92"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  +11.08%  "
$ws.Range("E9").Value = "  +3.07%  "
$ws.Range("E10").Value = "  +2.41%  "
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("D12").Value = "1.840.64"
$ws.Range("E12").Value = "  +3.65%  "
$ws.Range("D13").Value = "1.605.27"
$ws.Range("E13").Value = "  +3.27%  "
$ws.Range("D14").Value = "29.806.65"
$ws.Range("E14").Value = "  +4.86%  "
$ws.Range("E15").Value = "  +5.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.75"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  +3.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "244.39"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +6.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.24"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  +3.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.59"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +3.50%  "
$ws.Range("E20").Value = "  +3.11%  "
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.04"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  +3.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.23"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +3.60%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +3.76%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.21"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +2.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.32"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  +4.08%  "
$ws.Range("E27").Value = "  +5.23%  "
$ws.Range("E28").Value = "  +2.59%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0472"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  +1.24%  "
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.24"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +2.69%  "
$ws.Range("E33").Value = "  +4.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.11"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +3.48%  "
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.84"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +10.53%  "
$ws.Range("E37").Value = "  +2.12%  "
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("E39").Value = "  +2.95%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.534"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +4.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "55.27"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +28.37%  "
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.794"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.996"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0467"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.73"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +8.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.30"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("D48").Value = "1.751.20"
$ws.Range("E48").Value = "  +3.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.98"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +2.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.839"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -4.04%  "
